# "Fruta / hortaliza, semanal" — weekly data refresh for the
# "Feria Lagunitas de Puerto Montt - Cebollín" sheet.
#
# A new weekly price observation is inserted as row 148 (pushing the
# existing rows 148-185 down to 149-186 respectively); the new row
# shares the same market/category metadata as the rest of the block
# and carries its own date/volume/price figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 148, shifting the rest of the
# block (old 148..185) down to 149..186.
$ws.Rows.Item(148).Insert()

# Populate the newly inserted row 148 with the new observation.
$ws.Cells.Item(148, 1).Value = 4
$ws.Cells.Item(148, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(148, 3).Value = "Los Lagos"
$ws.Cells.Item(148, 4).Value = 44511
$ws.Cells.Item(148, 5).Value = 10
$ws.Cells.Item(148, 6).Value = 100112037
$ws.Cells.Item(148, 7).Value = "Cebollín"
$ws.Cells.Item(148, 8).Value = "Sin especificar"
$ws.Cells.Item(148, 9).Value = "Primera"
$ws.Cells.Item(148, 10).Value = 90
$ws.Cells.Item(148, 11).Value = 5000
$ws.Cells.Item(148, 12).Value = 5000
$ws.Cells.Item(148, 13).Value = 5000
$ws.Cells.Item(148, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(148, 15).Value = "Región Metropolitana"
$ws.Cells.Item(148, 16).Value = 139
$ws.Cells.Item(148, 17).Value = 36
$ws.Cells.Item(148, 18).Value = "Hortaliza"
